# Update the flight cost figures on the "Global" sheet with the new
# date/cost test data (commit: "new date and costs used gitbash couldnt
# do this in UFT").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

$ws.Range("C2").Value = 356.94
$ws.Range("C3").Value = 336.39999999999998
$ws.Range("C4").Value = 350.94
